$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-9 from 45175 (2023-09-06)
# to 45183 (2023-09-14), keeping existing formatting untouched.
$ws.Range("C2:C9").Value = 45183
